# Apply the "feat: add 2022-Q3 data" change.
#
# Before: 总计 | 2022-Q2 | 2022-Q1
# After:  总计 | 2022-Q3 (NEW) | 2022-Q2 | 2022-Q1
#
# 1. Insert a brand-new worksheet named "2022-Q3" right after "总计" and
#    fill it with the fund-holding table for that quarter.
# 2. Update the "总计" (summary) sheet so that it now lists three quarters
#    (2022-Q3, 2022-Q2, 2022-Q1) instead of two, with 2022-Q3's counts.
#
# NOTE: Worksheets.Item(N) appears to be a *positional* lookup rather
# than a stable object handle in this runtime - once a new sheet is
# inserted, a previously-captured Worksheets.Item(N) variable can end up
# referring to a different sheet than when it was captured. To avoid
# that pitfall, all sheet-insertion/renaming is done first, and any
# "by position" worksheet handles used afterwards are re-fetched fresh
# right before they're used.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: insert the new "2022-Q3" worksheet right after 总计 (position 1)
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)
$q3 = $wb.Worksheets.Add($null, $summary)
$q3.Name = "2022-Q3"
$q3.Outline.SummaryRow = 1
$q3.Outline.SummaryColumn = 1

# ---------------------------------------------------------------------
# Step 2: update 总计 (summary) sheet - still at position 1
# ---------------------------------------------------------------------

# Existing row 2 (was 2022-Q2 info) becomes the new 2022-Q3 row.
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 7
$summary.Range("D2").Value = 0.96

# Existing row 3 (was 2022-Q1 info) becomes the new 2022-Q2 row.
$summary.Range("B3").Value = "2022-Q2"
$summary.Range("C3").Value = 2
$summary.Range("D3").Value = 0

# Brand-new row 4 holds the 2022-Q1 row that used to be row 3.
# Copy the formatting of A3 (bold/centered/bordered) onto A4 first.
$summary.Range("A3").Copy()
$summary.Range("A4").PasteSpecial(-4122)
$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2022-Q1"
$summary.Range("C4").Value = 2
$summary.Range("D4").Value = 0.18

# ---------------------------------------------------------------------
# Step 3: fill in the "2022-Q3" sheet's data table
# ---------------------------------------------------------------------

# The "2022-Q2" sheet (now at position 3, after 总计 and the new
# 2022-Q3 sheet) already has a correctly-formatted (bold/centered/
# bordered) 7-column header row and count column - reuse it as the
# formatting template for the new sheet.
$refSheet = $wb.Worksheets.Item(3)

$refSheet.Range("B1:H1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)

$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# The B (fund code) and D/E/F/G columns in the source data are stored as
# text, not numbers (e.g. "688888", "013531", "8.15", "91.56") - notably
# some fund codes have leading zeros that must be preserved - so force
# those ranges to Text format before assigning the values.
$q3.Range("B2:B8").NumberFormat = "@"
$q3.Range("D2:G8").NumberFormat = "@"

$data = @(
    @(0, "688888", "浙商聚潮产业成长混合A", "8.15", "91.56", "4.58", "0.3733", 6),
    @(1, "013531", "浙商聚潮产业成长混合C", "5.49", "91.56", "4.58", "0.2514", 6),
    @(2, "010381", "浙商智选价值混合A",     "4.54", "91.07", "4.56", "0.2070", 7),
    @(3, "010382", "浙商智选价值混合C",     "2.31", "91.07", "4.56", "0.1053", 7),
    @(4, "013204", "恒生前海恒源天利债A",   "1.30", "32.00", "1.40", "0.0182", 7),
    @(5, "003670", "中融物联网主题灵活配置混合", "0.14", "85.56", "3.50", "0.0049", 5),
    @(6, "013205", "恒生前海恒源天利债C",   "0.00", "32.00", "1.40", $null, 7)
)

$rowIndex = 2
foreach ($entry in $data) {
    # Column A count cells use the bold/centered/bordered style too.
    $refSheet.Range("A2").Copy()
    $q3.Cells.Item($rowIndex, 1).PasteSpecial(-4122)

    $q3.Cells.Item($rowIndex, 1).Value = $entry[0]
    $q3.Cells.Item($rowIndex, 2).Value = $entry[1]
    $q3.Cells.Item($rowIndex, 3).Value = $entry[2]
    $q3.Cells.Item($rowIndex, 4).Value = $entry[3]
    $q3.Cells.Item($rowIndex, 5).Value = $entry[4]
    $q3.Cells.Item($rowIndex, 6).Value = $entry[5]
    $q3.Cells.Item($rowIndex, 7).Value = $entry[6]
    $q3.Cells.Item($rowIndex, 8).Value = $entry[7]
    $rowIndex = $rowIndex + 1
}

# Row 8's "持有市值(亿元)" (column G) is a true numeric 0, unlike the
# other rows in that column which are text. Clear the text formatting
# for that single cell and write it as a real number.
$q3.Range("G8").NumberFormat = "General"
$q3.Range("G8").Value = 0

# Keep "2022-Q1" (now the 4th sheet) as the selected/active tab, matching
# its original selected state.
$wb.Worksheets.Item(4).Activate()
